# B6-PowerPoint.pptx edit: 22 Mar 2020
#
# 1) Three tables (on the slides that hold them) switch from the deck's
#    custom "Table_0" style to the built-in PowerPoint table style
#    {551866AD-A619-40DD-9D89-7A4AFBEBFE55}.
# 2) The design theme is swapped: the deck had been using the "Integral"
#    (Red Violet) palette; it switches to the standard "Office" palette.

$p = $ppt.ActivePresentation

# --- 1) Re-style the three tables -----------------------------------------
$targetStyleId = "{551866AD-A619-40DD-9D89-7A4AFBEBFE55}"

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($targetStyleId)
        }
    }
}

# --- 2) Swap the theme colours from "Integral" (Red Violet) to "Office" ---
$colorScheme = $p.SlideMaster.ColorScheme

# Office theme palette, in the standard dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink
# order used by ColorScheme.Colors(1..12).
$officeColors = @(
    0,          # dk1      000000
    16777215,   # lt1      FFFFFF
    6968388,    # dk2      44546A
    15132391,   # lt2      E7E6E6
    13998939,   # accent1  5B9BD5
    3243501,    # accent2  ED7D31
    10855845,   # accent3  A5A5A5
    49407,      # accent4  FFC000
    12874308,   # accent5  4472C4
    4697456,    # accent6  70AD47
    12673797,   # hlink    0563C1
    7491477     # folHlink 954F72
)

for ($k = 1; $k -le 12; $k++) {
    $colorScheme.Colors($k).RGB = $officeColors[$k - 1]
}
